$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: write cell values for new rows 601-616, in row-major / column order ---
# row 601
$ws.Cells.Item(601,1).Value = 45190.8645544676
$ws.Cells.Item(601,2).Value = "mkdaniel2020@gmail.com"
$ws.Cells.Item(601,3).Value = "글로벌비즈니스"
$ws.Cells.Item(601,4).Value = 20226403
$ws.Cells.Item(601,5).Value = "권다넬"
$ws.Cells.Item(601,6).Value = "74:26"
$ws.Cells.Item(601,7).Value = 0.2
$ws.Cells.Item(601,8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(601,9).Value = "952만 명"
$ws.Cells.Item(601,10).Value = 0.059
$ws.Cells.Item(601,11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(601,12).Value = "Black"
$ws.Cells.Item(601,14).Value = "국민부담률을 아일랜드 수준으로 낮춰야 한다"

# row 602
$ws.Cells.Item(602,1).Value = 45190.87335079861
$ws.Cells.Item(602,2).Value = "clara9398@naver.com"
$ws.Cells.Item(602,3).Value = "미디어스쿨"
$ws.Cells.Item(602,4).Value = 20232552
$ws.Cells.Item(602,5).Value = "양은지"
$ws.Cells.Item(602,6).Value = "78:22"
$ws.Cells.Item(602,7).Value = 0.15
$ws.Cells.Item(602,8).Value = "우리나라의 국민부담률은 2010년 22.4%에서 꾸준히 상승하여 2020년 27.9%에 달하였다."
$ws.Cells.Item(602,9).Value = "779만 명"
$ws.Cells.Item(602,10).Value = 0.151
$ws.Cells.Item(602,11).Value = "중소기업이 신고법인수의 91%를 차지하는 데 부담하는 세액은 24.6%이다"
$ws.Cells.Item(602,12).Value = "Black"
$ws.Cells.Item(602,14).Value = "국민부담률을 아일랜드 수준으로 낮춰야 한다"

# row 603
$ws.Cells.Item(603,1).Value = 45190.9127390162
$ws.Cells.Item(603,2).Value = "plzmxn@naver.com"
$ws.Cells.Item(603,3).Value = "광고홍보학과"
$ws.Cells.Item(603,4).Value = 20232642
$ws.Cells.Item(603,5).Value = "현상희"
$ws.Cells.Item(603,6).Value = "75:25"
$ws.Cells.Item(603,7).Value = 0.2
$ws.Cells.Item(603,8).Value = "프랑스와 스웨덴의 국민부담률은 꾸준히 40%를 넘고 있다."
$ws.Cells.Item(603,9).Value = "952만 명"
$ws.Cells.Item(603,10).Value = 0.151
$ws.Cells.Item(603,11).Value = "상호출자제한기업은 신고법인수의 0.1%를 차지하는 데 부담하는 세액은 25.5%이다"
$ws.Cells.Item(603,12).Value = "Black"
$ws.Cells.Item(603,14).Value = "국민부담률을 OECD 평균 수준으로 높여야 한다"

# row 604
$ws.Cells.Item(604,1).Value = 45190.973911504625
$ws.Cells.Item(604,2).Value = "year0309@nate.com"
$ws.Cells.Item(604,3).Value = "경영학과"
$ws.Cells.Item(604,4).Value = 20222950
$ws.Cells.Item(604,5).Value = "박소희"
$ws.Cells.Item(604,6).Value = "75:25"
$ws.Cells.Item(604,7).Value = 0.2
$ws.Cells.Item(604,8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(604,9).Value = "952만 명"
$ws.Cells.Item(604,10).Value = 0.059
$ws.Cells.Item(604,11).Value = "중견기업은 신고법인수의 0.5%를 차지하는 데 부담하는 세액은 8.7%이다"
$ws.Cells.Item(604,12).Value = "Red"
$ws.Cells.Item(604,13).Value = "국민부담률을 OECD 평균 수준으로 높여야 한다"

# row 605
$ws.Cells.Item(605,1).Value = 45190.997704444446
$ws.Cells.Item(605,2).Value = "jgw1274@naver.com"
$ws.Cells.Item(605,3).Value = "화학과"
$ws.Cells.Item(605,4).Value = 20223427
$ws.Cells.Item(605,5).Value = "정근원"
$ws.Cells.Item(605,6).Value = "74:26"
$ws.Cells.Item(605,7).Value = 0.2
$ws.Cells.Item(605,8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(605,9).Value = "952만 명"
$ws.Cells.Item(605,10).Value = 0.059
$ws.Cells.Item(605,11).Value = "중견기업은 신고법인수의 0.5%를 차지하는 데 부담하는 세액은 8.7%이다"
$ws.Cells.Item(605,12).Value = "Black"
$ws.Cells.Item(605,14).Value = "국민부담률을 아일랜드 수준으로 낮춰야 한다"

# row 606
$ws.Cells.Item(606,1).Value = 45191.0694246412
$ws.Cells.Item(606,2).Value = "lsyun5050@naver.com"
$ws.Cells.Item(606,3).Value = "미디어스쿨"
$ws.Cells.Item(606,4).Value = 20232561
$ws.Cells.Item(606,5).Value = "이서연"
$ws.Cells.Item(606,6).Value = "74:26"
$ws.Cells.Item(606,7).Value = 0.2
$ws.Cells.Item(606,8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(606,9).Value = "952만 명"
$ws.Cells.Item(606,10).Value = 0.059
$ws.Cells.Item(606,11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(606,12).Value = "Black"
$ws.Cells.Item(606,14).Value = "국민부담률을 OECD 평균 수준으로 높여야 한다"

# row 607
$ws.Cells.Item(607,1).Value = 45191.13238966435
$ws.Cells.Item(607,2).Value = "aj4497@naver.com"
$ws.Cells.Item(607,3).Value = "반도체 디스플레이 스쿨"
$ws.Cells.Item(607,4).Value = 20233326
$ws.Cells.Item(607,5).Value = "안연지"
$ws.Cells.Item(607,6).Value = "77:23"
$ws.Cells.Item(607,7).Value = 0.2
$ws.Cells.Item(607,8).Value = "OECD평균은 2010년 31.6%에서 2020년 33.5%까지 상승하였다."
$ws.Cells.Item(607,9).Value = "166만 명"
$ws.Cells.Item(607,10).Value = 0.374
$ws.Cells.Item(607,11).Value = "그 외 대기업은 신고법인수의 8.3%를 차지하는 데 부담하는 세액은 41.2%이다"
$ws.Cells.Item(607,12).Value = "Black"
$ws.Cells.Item(607,14).Value = "모름/무응답"

# row 608
$ws.Cells.Item(608,1).Value = 45191.4775428125
$ws.Cells.Item(608,2).Value = "h20203411@glab.hallym.ac.kr"
$ws.Cells.Item(608,3).Value = "화학과"
$ws.Cells.Item(608,4).Value = 20203411
$ws.Cells.Item(608,5).Value = "반초원"
$ws.Cells.Item(608,6).Value = "74:26"
$ws.Cells.Item(608,7).Value = 0.2
$ws.Cells.Item(608,8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(608,9).Value = "952만 명"
$ws.Cells.Item(608,10).Value = 0.059
$ws.Cells.Item(608,11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(608,12).Value = "Black"
$ws.Cells.Item(608,14).Value = "국민부담률을 OECD 평균 수준으로 높여야 한다"

# row 609
$ws.Cells.Item(609,1).Value = 45191.533444189816
$ws.Cells.Item(609,2).Value = "dms_5236@naver.com"
$ws.Cells.Item(609,3).Value = "바이오메디컬학과"
$ws.Cells.Item(609,4).Value = 20233638
$ws.Cells.Item(609,5).Value = "장예은"
$ws.Cells.Item(609,6).Value = "76:24"
$ws.Cells.Item(609,7).Value = 0.2
$ws.Cells.Item(609,8).Value = "OECD평균은 2010년 31.6%에서 2020년 33.5%까지 상승하였다."
$ws.Cells.Item(609,9).Value = "166만 명"
$ws.Cells.Item(609,10).Value = 0.151
$ws.Cells.Item(609,11).Value = "상호출자제한기업은 신고법인수의 0.1%를 차지하는 데 부담하는 세액은 25.5%이다"
$ws.Cells.Item(609,12).Value = "Black"
$ws.Cells.Item(609,14).Value = "국민부담률을 OECD 평균 수준으로 높여야 한다"

# row 610
$ws.Cells.Item(610,1).Value = 45191.56543079861
$ws.Cells.Item(610,2).Value = "khjkhj0314@naver.com"
$ws.Cells.Item(610,3).Value = "소프트웨어학부"
$ws.Cells.Item(610,4).Value = 20235155
$ws.Cells.Item(610,5).Value = "김현종"
$ws.Cells.Item(610,6).Value = "77:23"
$ws.Cells.Item(610,7).Value = 0.2
$ws.Cells.Item(610,8).Value = "우리나라의 국민부담률은 2010년 22.4%에서 꾸준히 상승하여 2020년 27.9%에 달하였다."
$ws.Cells.Item(610,9).Value = "952만 명"
$ws.Cells.Item(610,10).Value = 0.059
$ws.Cells.Item(610,11).Value = "상호출자제한기업은 신고법인수의 0.1%를 차지하는 데 부담하는 세액은 25.5%이다"
$ws.Cells.Item(610,12).Value = "Black"
$ws.Cells.Item(610,14).Value = "모름/무응답"

# row 611
$ws.Cells.Item(611,1).Value = 45191.590850821754
$ws.Cells.Item(611,2).Value = "audri01@naver.com"
$ws.Cells.Item(611,3).Value = "바이오메디컬학과"
$ws.Cells.Item(611,4).Value = 20203639
$ws.Cells.Item(611,5).Value = "이승아"
$ws.Cells.Item(611,6).Value = "77:23"
$ws.Cells.Item(611,7).Value = 0.15
$ws.Cells.Item(611,8).Value = "우리나라의 국민부담률은 2010년 22.4%에서 꾸준히 상승하여 2020년 27.9%에 달하였다."
$ws.Cells.Item(611,9).Value = "166만 명"
$ws.Cells.Item(611,10).Value = 0.151
$ws.Cells.Item(611,11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(611,12).Value = "Black"
$ws.Cells.Item(611,14).Value = "모름/무응답"

# row 612
$ws.Cells.Item(612,1).Value = 45191.61094719908
$ws.Cells.Item(612,2).Value = "siwongim43@gmail.com"
$ws.Cells.Item(612,3).Value = "글로벌학부 정보법과학전공"
$ws.Cells.Item(612,4).Value = 20206410
$ws.Cells.Item(612,5).Value = "김시원"
$ws.Cells.Item(612,6).Value = "74:26"
$ws.Cells.Item(612,7).Value = 0.2
$ws.Cells.Item(612,8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(612,9).Value = "952만 명"
$ws.Cells.Item(612,10).Value = 0.059
$ws.Cells.Item(612,11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(612,12).Value = "Black"
$ws.Cells.Item(612,14).Value = "국민부담률을 OECD 평균 수준으로 높여야 한다"

# row 613
$ws.Cells.Item(613,1).Value = 45191.708008483794
$ws.Cells.Item(613,2).Value = "lih3488@naver.com"
$ws.Cells.Item(613,3).Value = "경영학과"
$ws.Cells.Item(613,4).Value = 20192953
$ws.Cells.Item(613,5).Value = "이인혁"
$ws.Cells.Item(613,6).Value = "77:23"
$ws.Cells.Item(613,7).Value = 0.15
$ws.Cells.Item(613,8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(613,9).Value = "166만 명"
$ws.Cells.Item(613,10).Value = 0.059
$ws.Cells.Item(613,11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(613,12).Value = "Black"
$ws.Cells.Item(613,14).Value = "국민부담률을 OECD 평균 수준으로 높여야 한다"

# row 614
$ws.Cells.Item(614,1).Value = 45191.78064038194
$ws.Cells.Item(614,2).Value = "limhyoeun8698@naver.com"
$ws.Cells.Item(614,3).Value = "미디어스쿨"
$ws.Cells.Item(614,4).Value = 20232572
$ws.Cells.Item(614,5).Value = "임효은"
$ws.Cells.Item(614,6).Value = "74:26"
$ws.Cells.Item(614,7).Value = 0.1
$ws.Cells.Item(614,8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(614,9).Value = "166만 명"
$ws.Cells.Item(614,10).Value = 0.151
$ws.Cells.Item(614,11).Value = "그 외 대기업은 신고법인수의 8.3%를 차지하는 데 부담하는 세액은 41.2%이다"
$ws.Cells.Item(614,12).Value = "Black"
$ws.Cells.Item(614,14).Value = "모름/무응답"

# row 615
$ws.Cells.Item(615,1).Value = 45191.78195478009
$ws.Cells.Item(615,2).Value = "zzolle424@naver.com"
$ws.Cells.Item(615,3).Value = "경영학과"
$ws.Cells.Item(615,4).Value = 20233030
$ws.Cells.Item(615,5).Value = "전지호"
$ws.Cells.Item(615,6).Value = "74:26"
$ws.Cells.Item(615,7).Value = 0.2
$ws.Cells.Item(615,8).Value = "OECD평균은 2010년 31.6%에서 2020년 33.5%까지 상승하였다."
$ws.Cells.Item(615,9).Value = "779만 명"
$ws.Cells.Item(615,10).Value = 0.151
$ws.Cells.Item(615,11).Value = "중소기업이 신고법인수의 91%를 차지하는 데 부담하는 세액은 24.6%이다"
$ws.Cells.Item(615,12).Value = "Red"
$ws.Cells.Item(615,13).Value = "국민부담률을 아일랜드 수준으로 낮춰야 한다"

# row 616
$ws.Cells.Item(616,1).Value = 45191.79859409722
$ws.Cells.Item(616,2).Value = "tlsfkdhs0321@naver.com"
$ws.Cells.Item(616,3).Value = "생명과학과"
$ws.Cells.Item(616,4).Value = 20233522
$ws.Cells.Item(616,5).Value = "신라온"
$ws.Cells.Item(616,6).Value = "74:26"
$ws.Cells.Item(616,7).Value = 0.2
$ws.Cells.Item(616,8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(616,9).Value = "952만 명"
$ws.Cells.Item(616,10).Value = 0.059
$ws.Cells.Item(616,11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(616,12).Value = "Black"
$ws.Cells.Item(616,14).Value = "모름/무응답"

# --- Step 2: apply per-column cell styles to match existing data rows ---
# Template row 594 (has A-L, N; no M) used for rows that populate column N
# Template row 599 (has A-L, M; no N) used for rows that populate column M

$src = $ws.Range("A594:L594")
$dst = $ws.Range("A601:L603")
$src.Copy()
$dst.PasteSpecial(-4122)
$src = $ws.Range("N594")
$dst = $ws.Range("N601:N603")
$src.Copy()
$dst.PasteSpecial(-4122)

$src = $ws.Range("A594:L594")
$dst = $ws.Range("A605:L614")
$src.Copy()
$dst.PasteSpecial(-4122)
$src = $ws.Range("N594")
$dst = $ws.Range("N605:N614")
$src.Copy()
$dst.PasteSpecial(-4122)

$src = $ws.Range("A594:L594")
$dst = $ws.Range("A616:L616")
$src.Copy()
$dst.PasteSpecial(-4122)
$src = $ws.Range("N594")
$dst = $ws.Range("N616:N616")
$src.Copy()
$dst.PasteSpecial(-4122)

$src = $ws.Range("A599:L599")
$dst = $ws.Range("A604:L604")
$src.Copy()
$dst.PasteSpecial(-4122)
$src = $ws.Range("M599")
$dst = $ws.Range("M604:M604")
$src.Copy()
$dst.PasteSpecial(-4122)

$src = $ws.Range("A599:L599")
$dst = $ws.Range("A615:L615")
$src.Copy()
$dst.PasteSpecial(-4122)
$src = $ws.Range("M599")
$dst = $ws.Range("M615:M615")
$src.Copy()
$dst.PasteSpecial(-4122)

$excel.CutCopyMode = 0
